$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.91848542527778
$ws.Range("C2").Value = 9.067267367789125
$ws.Range("D2").Value = 9.499106876410901
$ws.Range("E2").Value = 13.62590144479448
$ws.Range("F2").Value = 30.47522122678937
$ws.Range("I2").Value = 20.48835151856349
$ws.Range("J2").Value = 9.908396315229064
$ws.Range("M2").Value = 16.68047667615859
$ws.Range("O2").Value = 22.7119367478681
$ws.Range("B3").Value = 13.29949308189612
$ws.Range("C3").Value = 8.548197509756022
$ws.Range("D3").Value = 9.480893187953559
$ws.Range("E3").Value = 13.64908886173643
$ws.Range("F3").Value = 30.57476240681848
$ws.Range("I3").Value = 20.63619880020485
$ws.Range("J3").Value = 9.939428663107822
$ws.Range("M3").Value = 16.47699412924681
$ws.Range("O3").Value = 22.81734876917915
$ws.Range("B4").Value = 12.90461363996673
$ws.Range("C4").Value = 8.212467703286778
$ws.Range("D4").Value = 9.470867186237124
$ws.Range("E4").Value = 13.6655986273027
$ws.Range("F4").Value = 30.64500426764786
$ws.Range("I4").Value = 20.7322061525001
$ws.Range("J4").Value = 9.959712905702959
$ws.Range("M4").Value = 16.35258396588184
$ws.Range("O4").Value = 22.88849551488509
$ws.Range("B5").Value = 12.74017034032925
$ws.Range("C5").Value = 8.071460690999109
$ws.Range("D5").Value = 9.467075182654199
$ws.Range("E5").Value = 13.67289778988152
$ws.Range("F5").Value = 30.67591310906784
$ws.Range("I5").Value = 20.77264501360665
$ws.Range("J5").Value = 9.968288673560075
$ws.Range("M5").Value = 16.30206831379333
$ws.Range("O5").Value = 22.91909751714355
$ws.Range("B6").Value = 12.71265802541993
$ws.Range("C6").Value = 8.047796123265611
$ws.Range("D6").Value = 9.466463334914298
$ws.Range("E6").Value = 13.67414430885814
$ws.Range("F6").Value = 30.68118320475525
$ws.Range("I6").Value = 20.77943929512415
$ws.Range("J6").Value = 9.969731395240117
$ws.Range("M6").Value = 16.29369269332696
$ws.Range("O6").Value = 22.92427596427829
$ws.Range("B7").Value = 12.9024098959611
$ws.Range("C7").Value = 8.210582889769672
$ws.Range("D7").Value = 9.470814853481842
$ws.Range("E7").Value = 13.66569475347215
$ws.Range("F7").Value = 30.64541187721266
$ws.Range("I7").Value = 20.73274619990976
$ws.Range("J7").Value = 9.959827306688304
$ws.Range("M7").Value = 16.35190188996009
$ws.Range("O7").Value = 22.88890171702106
$ws.Range("B8").Value = 13.70824246506742
$ws.Range("C8").Value = 8.891869731625343
$ws.Range("D8").Value = 9.492588254213532
$ws.Range("E8").Value = 13.63342471333345
$ws.Range("F8").Value = 30.50764403209459
$ws.Range("I8").Value = 20.53824426593468
$ws.Range("J8").Value = 9.918841179649135
$ws.Range("M8").Value = 16.61023583591012
$ws.Range("O8").Value = 22.74694602450092
$ws.Range("B9").Value = 15.1635561742484
$ws.Range("C9").Value = 10.08996448237642
$ws.Range("D9").Value = 9.544344730787657
$ws.Range("E9").Value = 13.58818032652293
$ws.Range("F9").Value = 30.31026547139905
$ws.Range("I9").Value = 20.19829789033787
$ws.Range("J9").Value = 9.848210174333461
$ws.Range("M9").Value = 17.1188522632496
$ws.Range("O9").Value = 22.51980530862416
$ws.Range("B10").Value = 16.14844179053456
$ws.Range("C10").Value = 10.88326654994804
$ws.Range("D10").Value = 9.587717946248635
$ws.Range("E10").Value = 13.56593848627274
$ws.Range("F10").Value = 30.21012615203031
$ws.Range("I10").Value = 19.97380191811712
$ws.Range("J10").Value = 9.802230790703762
$ws.Range("M10").Value = 17.49087948684345
$ws.Range("O10").Value = 22.38451855140818
$ws.Range("B11").Value = 16.57686778749527
$ws.Range("C11").Value = 11.22489590935839
$ws.Range("D11").Value = 9.608570010276937
$ws.Range("E11").Value = 13.55820732634665
$ws.Range("F11").Value = 30.17440855321444
$ws.Range("I11").Value = 19.87715656289764
$ws.Range("J11").Value = 9.782592147088407
$ws.Range("M11").Value = 17.65916527553283
$ws.Range("O11").Value = 22.32990817964417
$ws.Range("B12").Value = 16.73619600131225
$ws.Range("C12").Value = 11.35147220747716
$ws.Range("D12").Value = 9.616623646565142
$ws.Range("E12").Value = 13.55562265720552
$ws.Range("F12").Value = 30.16230423725458
$ws.Range("I12").Value = 19.84134766103173
$ws.Range("J12").Value = 9.775338837415552
$ws.Range("M12").Value = 17.72270502882546
$ws.Range("O12").Value = 22.31023149954709
$ws.Range("B13").Value = 16.70201238164894
$ws.Range("C13").Value = 11.32433618867892
$ws.Range("D13").Value = 9.614882213777477
$ws.Range("E13").Value = 13.55616406432212
$ws.Range("F13").Value = 30.164847813023
$ws.Range("I13").Value = 19.84902466271054
$ws.Range("J13").Value = 9.776892814530925
$ws.Range("M13").Value = 17.70902967325288
$ws.Range("O13").Value = 22.31442451476305
$ws.Range("B14").Value = 16.59003453921537
$ws.Range("C14").Value = 11.2353654718625
$ws.Range("D14").Value = 9.609229455566016
$ws.Range("E14").Value = 13.55798781264551
$ws.Range("F14").Value = 30.17338420593044
$ws.Range("I14").Value = 19.87419473370509
$ws.Range("J14").Value = 9.781991738644876
$ws.Range("M14").Value = 17.66439674499694
$ws.Range("O14").Value = 22.32826922893281
$ws.Range("B15").Value = 16.52106387158308
$ws.Range("C15").Value = 11.18050423790627
$ws.Range("D15").Value = 9.605787360972132
$ws.Range("E15").Value = 13.5591495642071
$ws.Range("F15").Value = 30.17879825653215
$ws.Range("I15").Value = 19.88971485134631
$ws.Range("J15").Value = 9.785138854634551
$ws.Range("M15").Value = 17.63703205963807
$ws.Range("O15").Value = 22.33688032932679
$ws.Range("B16").Value = 16.12004076985615
$ws.Range("C16").Value = 10.86055079218005
$ws.Range("D16").Value = 9.586377455139816
$ws.Range("E16").Value = 13.56649174377899
$ws.Range("F16").Value = 30.21265890430877
$ws.Range("I16").Value = 19.98022822950673
$ws.Range("J16").Value = 9.803539905143023
$ws.Range("M16").Value = 17.47985820052676
$ws.Range("O16").Value = 22.3882274986926
$ws.Range("B17").Value = 15.86893742706927
$ws.Range("C17").Value = 10.65932102884788
$ws.Range("D17").Value = 9.574754636227246
$ws.Range("E17").Value = 13.5716070899467
$ws.Range("F17").Value = 30.23595564852512
$ws.Range("I17").Value = 20.0371590288902
$ws.Range("J17").Value = 9.815155363474606
$ws.Range("M17").Value = 17.38315776848793
$ws.Range("O17").Value = 22.42150732452356
$ws.Range("B18").Value = 15.72266852319923
$ws.Range("C18").Value = 10.54176824124067
$ws.Range("D18").Value = 9.568175287156018
$ws.Range("E18").Value = 13.57477397778677
$ws.Range("F18").Value = 30.25028080320361
$ws.Range("I18").Value = 20.07041977551376
$ws.Range("J18").Value = 9.821956539827131
$ws.Range("M18").Value = 17.32745159362586
$ws.Range("O18").Value = 22.4413011331992
$ws.Range("B19").Value = 15.67283115085608
$ws.Range("C19").Value = 10.5016567153999
$ws.Range("D19").Value = 9.565965916297506
$ws.Range("E19").Value = 13.57588482607767
$ws.Range("F19").Value = 30.25528978974972
$ws.Range("I19").Value = 20.08176983749194
$ws.Range("J19").Value = 9.824279967397091
$ws.Range("M19").Value = 17.30857709222536
$ws.Range("O19").Value = 22.44811478875508
$ws.Range("B20").Value = 15.89585902090505
$ws.Range("C20").Value = 10.68092987761986
$ws.Range("D20").Value = 9.575980983690197
$ws.Range("E20").Value = 13.57103930103277
$ws.Range("F20").Value = 30.23337983508215
$ws.Range("I20").Value = 20.03104527473502
$ws.Range("J20").Value = 9.813906432610658
$ws.Range("M20").Value = 17.39346102556127
$ws.Range("O20").Value = 22.41789709124472
$ws.Range("B21").Value = 16.62300468485926
$ws.Range("C21").Value = 11.26157423952492
$ws.Range("D21").Value = 9.61088556694717
$ws.Range("E21").Value = 13.5574428285175
$ws.Range("F21").Value = 30.17083824078531
$ws.Range("I21").Value = 19.86678026498458
$ws.Range("J21").Value = 9.780489085456232
$ws.Range("M21").Value = 17.67751196631068
$ws.Range("O21").Value = 22.32417542852955
$ws.Range("B22").Value = 17.08125866131065
$ws.Range("C22").Value = 11.62478229008982
$ws.Range("D22").Value = 9.634613356513858
$ws.Range("E22").Value = 13.55055554384577
$ws.Range("F22").Value = 30.13824987805274
$ws.Range("I22").Value = 19.76402056382521
$ws.Range("J22").Value = 9.759717887474354
$ws.Range("M22").Value = 17.86204702248827
$ws.Range("O22").Value = 22.26877241194505
$ws.Range("B23").Value = 16.83825889104382
$ws.Range("C23").Value = 11.43242679329301
$ws.Range("D23").Value = 9.621866918084615
$ws.Range("E23").Value = 13.55404864169849
$ws.Range("F23").Value = 30.15488269889446
$ws.Range("I23").Value = 19.81844445159152
$ws.Range("J23").Value = 9.770706152753457
$ws.Range("M23").Value = 17.76367467475554
$ws.Range("O23").Value = 22.29780479713165
$ws.Range("B24").Value = 15.8836936990162
$ws.Range("C24").Value = 10.67116631848271
$ws.Range("D24").Value = 9.575426231624924
$ws.Range("E24").Value = 13.57129529439372
$ws.Range("F24").Value = 30.23454145912152
$ws.Range("I24").Value = 20.03380765060998
$ws.Range("J24").Value = 9.814470690172159
$ws.Range("M24").Value = 17.38880326749988
$ws.Range("O24").Value = 22.41952721981652
$ws.Range("B25").Value = 14.78415811793451
$ws.Range("C25").Value = 9.780953717132689
$ws.Range("D25").Value = 9.52939075391466
$ws.Range("E25").Value = 13.59848840861693
$ws.Range("F25").Value = 30.35581501023114
$ws.Range("I25").Value = 20.28582505695549
$ws.Range("J25").Value = 9.866277451939727
$ws.Range("M25").Value = 16.98134152761173
$ws.Range("O25").Value = 22.57573033383644
